# Insert a new weekly data row for "Cebollín" (Primera, Vega Monumental
# Concepción) at row 32. Inserting shifts every existing row from 32..57
# down to 33..58 (and grows the sheet's used range to A1:R58), matching
# the target diff exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(32).Insert()

$ws.Range("A32").Value = 11
$ws.Range("B32").Value = "Vega Monumental Concepción"
$ws.Range("C32").Value = "Bíobío"
$ws.Range("D32").Value = 44763
$ws.Range("E32").Value = 8
$ws.Range("F32").Value = 100112037
$ws.Range("G32").Value = "Cebollín"
$ws.Range("H32").Value = "Sin especificar"
$ws.Range("I32").Value = "Primera"
$ws.Range("J32").Value = 80
$ws.Range("K32").Value = 5000
$ws.Range("L32").Value = 5500
$ws.Range("M32").Value = 5188
$ws.Range("N32").Value = "`$/paquete 36 unidades"
$ws.Range("O32").Value = "Región Metropolitana"
$ws.Range("P32").Value = 144
$ws.Range("Q32").Value = 36
$ws.Range("R32").Value = "Hortaliza"
